# Insert a new data row at row 87 (a new daily price observation for
# Navel Late oranges), shifting the existing rows 87-187 down to 88-188.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("87:87").Insert()

$ws.Range("A87").Value = 11
$ws.Range("B87").Value = "Vega Monumental Concepción"
$ws.Range("C87").Value = "Bíobío"
$ws.Range("D87").Value = 44539
$ws.Range("E87").Value = 8
$ws.Range("F87").Value = "Fruta"
$ws.Range("G87").Value = 100102
$ws.Range("H87").Value = "Cítricos"
$ws.Range("I87").Value = 100102005
$ws.Range("J87").Value = "Naranja"
$ws.Range("K87").Value = "Navel Late"
$ws.Range("L87").Value = "Primera"
$ws.Range("M87").Value = 220
$ws.Range("N87").Value = 7500
$ws.Range("O87").Value = 8000
$ws.Range("P87").Value = 7727
$ws.Range("Q87").Value = "$/caja 15 kilos granel"
$ws.Range("R87").Value = "Región de O'Higgins"
$ws.Range("S87").Value = 515
$ws.Range("T87").Value = 15
